$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (coinranking snapshot).
# Price cells that look like plain decimals are written with a leading
# apostrophe so Excel stores them as text (matching the existing text-based
# Price/Volume columns) instead of auto-converting them to numbers.

$ws.Range("D2").Value = "27.382.70"
$ws.Range("D3").Value = "1.654.93"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'213.17"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "'0.539"
$ws.Range("E6").Value = "  +5.22%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'23.47"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").Value = "'0.261"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").Value = "'0.0906"
$ws.Range("E11").Value = "  +3.40%  "
$ws.Range("D12").Value = "1.889.27"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "1.639.20"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "'0.567"
$ws.Range("E15").Value = "  +3.69%  "
$ws.Range("D16").Value = "'65.53"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "27.386.81"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "'229.36"
$ws.Range("E18").Value = "  -6.96%  "
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").Value = "'7.43"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  -2.72%  "
$ws.Range("D23").Value = "'9.40"
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("D24").Value = "'2.04"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").Value = "'147.02"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'7.06"
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.115"
$ws.Range("E27").Value = "  +3.01%  "
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'15.71"
$ws.Range("E29").Value = "  -3.03%  "
$ws.Range("D30").Value = "'0.0494"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("E31").Value = "  -3.89%  "
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").Value = "1.422.46"
$ws.Range("E34").Value = "  -1.75%  "
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("E37").Value = "  -2.60%  "
$ws.Range("D38").Value = "'0.568"
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "'1.04"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  +2.78%  "
$ws.Range("D43").Value = "'65.03"
$ws.Range("E43").Value = "  -5.69%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'0.789"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.797.60"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").Value = "'1.66"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "'88.03"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  -3.51%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "'7.71"
$ws.Range("E51").Value = "  -1.29%  "
